# "Screenshot Code added to TC06" — populate newly-captured automation
# evidence (employee/user ids + names pulled from the latest test run)
# into the TestData workbook's per-module sheets.

$wb = $excel.ActiveWorkbook

# --- AddEmployee -----------------------------------------------------
$ws = $wb.Worksheets.Item("AddEmployee")
$ws.Range("F2").Value = "Val"
$ws.Range("G2").Value = "Johnson"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "0150"
$ws.Range("F3").Value = "Dalton"
$ws.Range("G3").Value = "Tillman"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "0151"
$ws.Range("F5").Value = "Donnie"
$ws.Range("G5").Value = "Legros"
$ws.Range("F6").Value = "Gregory"
$ws.Range("G6").Value = "Greenholt"

# --- AddUser -----------------------------------------------------------
$ws = $wb.Worksheets.Item("AddUser")
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "1110"
$ws.Range("G3").Value = "reiko.bradtke"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "1057"
$ws.Range("G6").Value = "lesia.kiehn"

# --- EditEmployee (TC06) ------------------------------------------------
$ws = $wb.Worksheets.Item("EditEmployee")
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "1075"
$ws.Range("G3").Value = "Leuschke"
$ws.Range("I3").Value = "India Office"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "0146"
$ws.Range("G6").Value = "Harber"
$ws.Range("I3").Select()

# --- RecruitmentCandidate ------------------------------------------------
$ws = $wb.Worksheets.Item("RecruitmentCandidate")
$ws.Range("F6").Value = "Paul"
$ws.Range("G6").Value = "Hills"

# --- AllOrdersTotal: move the active selection only ---------------------
$ws = $wb.Worksheets.Item("AllOrdersTotal")
$ws.Range("F6").Select()
